# Insert a new "DemonstrationProjectIdentifier" column right before the
# existing "Note" column (column DL), shifting every column from DL
# onward (both the header row and the data row) one position to the
# right. This matches the dimension change from A1:MN2 to A1:MO2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("DL").Insert()
$ws.Range("DL1").Value = "DemonstrationProjectIdentifier"

# Update the row-2 record id to the new value from the commit.
$ws.Range("A2").Value = "690148897e79911955eafc50"
